# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    # Force the literal string into the cell as TEXT (not an auto-parsed
    # number/percentage) without leaving a stray number-format style behind:
    # flip to text format, assign, then restore General/Normal so the
    # cell's style index matches its original (unstyled) state.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "65.217.97"
Set-TextCell $ws "E2" "  -5.93%  "
Set-TextCell $ws "D3" "3.451.21"
Set-TextCell $ws "E3" "  -7.40%  "
Set-TextCell $ws "D4" "1.00"
Set-TextCell $ws "E4" "  +0.24%  "
Set-TextCell $ws "D5" "554.80"
Set-TextCell $ws "E5" "  -9.46%  "
Set-TextCell $ws "D6" "180.50"
Set-TextCell $ws "E6" "  -6.15%  "
Set-TextCell $ws "D7" "3.449.26"
Set-TextCell $ws "E7" "  -7.29%  "
Set-TextCell $ws "D8" "0.595"
Set-TextCell $ws "E8" "  -6.75%  "
Set-TextCell $ws "D9" "1.00"
Set-TextCell $ws "E9" "  -0.03%  "
Set-TextCell $ws "D10" "0.640"
Set-TextCell $ws "E10" "  -12.20%  "
Set-TextCell $ws "D11" "0.140"
Set-TextCell $ws "E11" "  -13.76%  "
Set-TextCell $ws "D12" "50.86"
Set-TextCell $ws "E12" "  -15.89%  "
Set-TextCell $ws "D13" "0.0000248"
Set-TextCell $ws "E13" "  -14.68%  "
Set-TextCell $ws "D14" "9.36"
Set-TextCell $ws "E14" "  -12.45%  "
Set-TextCell $ws "D15" "4.002.96"
Set-TextCell $ws "E15" "  -7.27%  "
Set-TextCell $ws "E16" "  -1.74%  "
Set-TextCell $ws "D17" "3.450.40"
Set-TextCell $ws "E17" "  -7.39%  "
Set-TextCell $ws "D18" "64.973.53"
Set-TextCell $ws "E18" "  -6.02%  "
Set-TextCell $ws "D19" "17.58"
Set-TextCell $ws "E19" "  -9.67%  "
Set-TextCell $ws "D20" "11.54"
Set-TextCell $ws "E20" "  -10.89%  "
Set-TextCell $ws "D21" "1.02"
Set-TextCell $ws "E21" "  -11.26%  "
Set-TextCell $ws "D22" "373.60"
Set-TextCell $ws "D23" "4.06"
Set-TextCell $ws "E23" "  -11.07%  "
Set-TextCell $ws "D24" "81.99"
Set-TextCell $ws "E24" "  -8.67%  "
Set-TextCell $ws "D25" "10.65"
Set-TextCell $ws "E25" "  -2.28%  "
Set-TextCell $ws "D26" "5.99"
Set-TextCell $ws "E26" "  -0.97%  "
Set-TextCell $ws "D27" "2.76"
Set-TextCell $ws "E27" "  -9.98%  "
Set-TextCell $ws "D28" "11.72"
Set-TextCell $ws "E28" "  -9.04%  "
Set-TextCell $ws "D29" "3.37"
Set-TextCell $ws "E29" "  -11.48%  "
Set-TextCell $ws "D30" "8.43"
Set-TextCell $ws "E30" "  -13.02%  "
Set-TextCell $ws "D31" "29.87"
Set-TextCell $ws "E31" "  -9.86%  "
Set-TextCell $ws "D32" "7.06"
Set-TextCell $ws "E32" "  -8.73%  "
Set-TextCell $ws "D33" "608.92"
Set-TextCell $ws "E33" "  -4.30%  "
Set-TextCell $ws "D34" "11.72"
Set-TextCell $ws "E34" "  -8.35%  "
Set-TextCell $ws "D35" "62.32"
Set-TextCell $ws "E35" "  -5.79%  "
Set-TextCell $ws "E36" "  -12.86%  "
Set-TextCell $ws "D37" "40.07"
Set-TextCell $ws "E37" "  -12.52%  "
Set-TextCell $ws "D38" "0.999"
Set-TextCell $ws "E38" "  -0.10%  "
Set-TextCell $ws "D39" "0.389"
Set-TextCell $ws "E39" "  -6.84%  "
Set-TextCell $ws "D40" "0.997"
Set-TextCell $ws "E40" "  -0.32%  "
Set-TextCell $ws "D41" "0.0₃0707"
Set-TextCell $ws "E41" "  -14.85%  "
Set-TextCell $ws "D42" "0.128"
Set-TextCell $ws "E42" "  -9.36%  "
Set-TextCell $ws "D43" "2.900.63"
Set-TextCell $ws "E43" "  +0.60%  "
Set-TextCell $ws "D44" "2.68"
Set-TextCell $ws "E44" "  -12.44%  "
Set-TextCell $ws "D45" "2.41"
Set-TextCell $ws "E45" "  -8.23%  "
Set-TextCell $ws "D46" "3.10"
Set-TextCell $ws "E46" "  -0.02%  "
Set-TextCell $ws "D47" "0.0389"
Set-TextCell $ws "E47" "  -13.44%  "
Set-TextCell $ws "D48" "0.126"
Set-TextCell $ws "E48" "  -10.47%  "
Set-TextCell $ws "D49" "137.36"
Set-TextCell $ws "E49" "  -3.36%  "
Set-TextCell $ws "D50" "2.42"
Set-TextCell $ws "E50" "  -11.78%  "
Set-TextCell $ws "D51" "8.04"
Set-TextCell $ws "E51" "  -12.23%  "
